# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计"),
#    populated with fund-holding detail rows.
# 2. Insert a new first data row into "总计" summarising the new quarter,
#    pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: new "2022-Q1" worksheet
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Add($null, $anchor)
$ws.Name = "2022-Q1"

# Pull the header / index-column styling from an existing sheet (style
# index 2 in the shared styles table: bold, bordered, centered) so we
# reuse it instead of synthesising a brand-new cell format.
$styleSrc = $wb.Worksheets.Item("2021-Q4")

$styleSrc.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Text-valued columns that look numeric (code / scale / weight figures)
# must stay text, matching the source data ("@" = text number format).
$ws.Range("B2:B11").NumberFormat = "@"
$ws.Range("D2:G11").NumberFormat = "@"

$rows = @(
    @(0, "159865", "国泰中证畜牧养殖ETF", "24.32", "99.29", "3.46", "0.8415", 9),
    @(1, "460007", "华泰柏瑞行业领先混合", "3.45", "94.37", "6.66", "0.2298", 6),
    @(2, "159867", "鹏华中证畜牧养殖ETF", "5.61", "97.87", "3.58", "0.2008", 9),
    @(3, "001398", "华泰柏瑞健康生活灵活配置混合", "3.09", "94.36", "6.30", "0.1947", 6),
    @(4, "011111", "华泰柏瑞行业严选混合型证券投资基金A", "2.36", "94.67", "7.03", "0.1659", 7),
    @(5, "516760", "平安中证畜牧养殖ETF", "1.45", "97.82", "3.58", "0.0519", 9),
    @(6, "516670", "招商中证畜牧养殖ETF", "1.06", "98.73", "3.61", "0.0383", 9),
    @(7, "011112", "华泰柏瑞行业严选混合型证券投资基金C", "0.26", "94.67", "7.03", "0.0183", 7),
    @(8, "000417", "国联安新精选灵活配置混合", "0.52", "28.16", "2.44", "0.0127", 2),
    @(9, "003981", "中银证券瑞益灵活配置混合C", "0.21", "89.21", "5.50", "0.0116", 2)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $styleSrc.Range("A2").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("A$r").Value = $data[0]
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("C$r").Value = $data[2]
    $ws.Range("D$r").Value = $data[3]
    $ws.Range("E$r").Value = $data[4]
    $ws.Range("F$r").Value = $data[5]
    $ws.Range("G$r").Value = $data[6]
    $ws.Range("H$r").Value = $data[7]
}

# ---------------------------------------------------------------------
# Part 2: prepend a "2022-Q1" summary row into "总计"
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Snapshot the current data rows (2..6) before shifting them down.
$totVals = @()
for ($r = 2; $r -le 6; $r++) {
    $row = @($tot.Range("A$r").Value(), $tot.Range("B$r").Value(), $tot.Range("C$r").Value(), $tot.Range("D$r").Value())
    $totVals += ,$row
}

# Write them back one row lower (bottom-up so we never clobber a row
# before it has been read), carrying the index-column style (2) along
# via a single-cell format copy rather than Rows.Insert (which would
# synthesise a brand-new, unused style entry).
for ($i = $totVals.Length - 1; $i -ge 0; $i--) {
    $src = $i + 2
    $dst = $i + 3

    $tot.Range("A$src").Copy()
    $tot.Range("A$dst").PasteSpecial(-4122)

    $tot.Range("A$dst").Value = $totVals[$i][0]
    $tot.Range("B$dst").Value = $totVals[$i][1]
    $tot.Range("C$dst").Value = $totVals[$i][2]
    $tot.Range("D$dst").Value = $totVals[$i][3]
}

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 10
$tot.Range("D2").Value = 1.77
